# "Added module 9 - PWM" - mark attendance for module dates 09/03/2021 (col D)
# and 16/03/2021 (col E) for the "2B" class, while filtering both the
# "Presenças" and "Material" tables down to class "2B".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Presenças "
$ws2 = $wb.Worksheets.Item(2)   # "Material"

# Filter both tables on the "Turma" column (3rd column) to show only class "2B".
$table1 = $ws1.ListObjects.Item(1)
$table1.Range.AutoFilter(3, @("2B"), 7)

$table2 = $ws2.ListObjects.Item(1)
$table2.Range.AutoFilter(3, @("2B"), 7)

# Mark attendance (value 1) in columns D (09/03/2021) and E (16/03/2021)
# for every student in class "2B" (the rows left visible by the filter).
$deRows = @(3,6,7,8,16,17,18,19,20,22,26,27,28,29)
foreach ($r in $deRows) {
    $ws1.Range("D$r").Value = 1
    $ws1.Range("E$r").Value = 1
}

# An extra (empty) underlined cell was left below the table.
$ws1.Range("E30").Font.Underline = 2

# Restore the selections / active sheet seen in the saved workbook.
$ws2.Range("F6").Select()
$ws1.Activate()
$ws1.Range("E30").Select()
